$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
# Target canonical widths (chars): B=16.9, C=11.7, E=13
# The engine quantizes ColumnWidth assignment to 1/6-character pixel steps,
# so we feed the input value whose quantized result is closest to the target.
$ws.Columns.Item(2).ColumnWidth = 16
$ws.Columns.Item(3).ColumnWidth = 10.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.166666666666666

# --- Row 2 data changes ---
$ws.Range("A2").Value = "CASE11591"
$ws.Range("B2").Value = "Authorization"
$ws.Range("C2").Value = "Completed"
$ws.Range("E2").Value = "TEST USER"
$ws.Range("F2").Value = "05-Jun-2024 06:46 PM"
